# Enemy and Summon AI doc update
$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Section 1: after "Once enemy is seen, Ghoul will charge at enemy."
# add a trailing space, then two new ilvl=2 bullet paragraphs.
# ---------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4.Range.InsertAfter(" ")

$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "The Ghoul is a ranged fighter with his poisonous snot balls but can switch into melee if need be with its paralyzing touch."
$p5.Range.ListFormat.ListLevelNumber = 3

$p5 = $d.Paragraphs(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "If a ghoul sees another enemy kill it becomes fearful and runs away."
$p6.Range.ListFormat.ListLevelNumber = 3

# ---------------------------------------------------------------
# Section 2: before the final "Summon" paragraph (which carries the
# _GoBack bookmark), insert the Summon/Fat Abomination/Babe Warrior
# outline block. Insert each new bullet immediately after the
# previous one so the block stays in document order.
# ---------------------------------------------------------------
function Insert-BulletAfter($afterIndex, $text, $listLevelNumber) {
    $anchor = $d.Paragraphs($afterIndex)
    $anchor.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs($newIndex)
    $newPara.Range.Text = $text
    $newPara.Range.ListFormat.ListLevelNumber = $listLevelNumber
    return $newIndex
}

# Index of the paragraph that currently ends with "...Continue air melee."
# (the anchor right before the final "Summon" + bookmark paragraph).
$cursor = $d.Paragraphs.Count - 1

$cursor = Insert-BulletAfter $cursor "Summon" 2
$cursor = Insert-BulletAfter $cursor "The summoned ghoul should attack the same way as the enemy ghoul, only he runs away if another summon is killed." 3
$cursor = Insert-BulletAfter $cursor "Fat Abomination" 1
$cursor = Insert-BulletAfter $cursor "Fat Abomination is a slow moving fighter, with a few tricks up his sleeves." 2
$cursor = Insert-BulletAfter $cursor "He can jump pretty high and slam down on enemies that are beneath him." 2
$cursor = Insert-BulletAfter $cursor "He can also grab at ranged targets with his hook and drag them to him…then slash them with hooks." 2
$cursor = Insert-BulletAfter $cursor "Fat Abomination walks slowly toward target if seen." 3
$cursor = Insert-BulletAfter $cursor "If target is range and target is in hooks path and Hook throw is not on a CD… Fat Abomination will throw hook." 4
$cursor = Insert-BulletAfter $cursor "If target is hooked and is now close to Fat Abomination, he will attack at melee with hooks." 5
$cursor = Insert-BulletAfter $cursor "If target is not hook…Fat Abomination will move toward target." 5
$cursor = Insert-BulletAfter $cursor "If target is in range and target is in hooks path and Hook throw is on a CD…Fat Abomination will continue to walk toward target." 4
$cursor = Insert-BulletAfter $cursor "If target is in range of Jump Slam and Jump Slam is not on CD…Fat Abomination will Jump and do Jump Slam." 4
$cursor = Insert-BulletAfter $cursor "If targets is in range of Jump Slam and Hook Grab…and Jump Slam is not on CD and Hook Grab is not on CD…Fat Abom will jump and do a jump slam." 4
$cursor = Insert-BulletAfter $cursor "If targets is in range of Jump Slam and Hook Grab…and Jump Slam is on on CD and Hook Grab is on CD…Fat Abom will melee closes enemy." 4
$cursor = Insert-BulletAfter $cursor "If target is in melee range Fat Abomination will do a melee attack." 4
$cursor = Insert-BulletAfter $cursor "Babe Warrior" 1
$cursor = Insert-BulletAfter $cursor "Enemy" 2

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
